$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of column S (the existing "2022" column) into the
# new column T so the added "2023" column matches the workbook's existing
# per-row number formats/fonts/borders.
$ws.Range("S4:S14").Copy()
$ws.Range("T4:T14").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header for the new year column
$ws.Range("T4").Value = 2023

# Data values for the new "2023" column (one per indicator row)
$ws.Range("T5").Value = 4.8187602774004432
$ws.Range("T6").Value = 11.788953009068425
$ws.Range("T7").Value = 5.2855407047387608
$ws.Range("T8").Value = 11.35112240576027
$ws.Range("T9").Value = 16.577540106951872
$ws.Range("T10").Value = 14.651002073255009
$ws.Range("T11").Value = 5.034965034965035
$ws.Range("T12").Value = 3.1837160751565765
$ws.Range("T13").Value = 2.2263731825525039
$ws.Range("T14").Value = 5.1321450522433931
